$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (date style) from A2 down to A3 so the new date cell
# reuses the existing cell style instead of creating a new one.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# New trade/result row (row 3) - same layout as row 2, new values
$ws.Range("A3").Value = 42605.885254629633
$ws.Range("B3").Value = 17
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "Random"
